# Auto commit at 2025-11-06 12:17:41.39
# Updates the metric values on the "Metrics" sheet; downstream formulas on
# the "today" sheet (which reference Metrics!B2:B13) recalculate
# automatically. Also restores the active-cell selection recorded in each
# sheet's view at save time.

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value = 65400.319999999992
$metrics.Range("B3").Value = 57101.499999999993
$metrics.Range("B4").Value = 20223.8
$metrics.Range("B5").Value = 2741
$metrics.Range("B6").Value = 4861646.0699999994
$metrics.Range("B7").Value = 4099178.1800000006
$metrics.Range("B8").Value = 1427183.63
$metrics.Range("B9").Value = 188948
$metrics.Range("B10").Value = 33327027.060000002
$metrics.Range("B11").Value = 31374453.34
$metrics.Range("B12").Value = 11708905.670000002
$metrics.Range("B13").Value = 1286578

# Selection on "Metrics" moves to E20 (single cell, no multi-cell sqref).
$metrics.Range("E20").Select()

$today = $wb.Worksheets.Item("today")

# Selection on "today" moves to G6 (single cell).
$today.Select()
$today.Range("G6").Select()
